$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original plain-text storage
# (values like "1.00" or "501.79" must not be auto-converted to numbers).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "54.537.17"
$ws.Range("E2").Value = "  +0.20%  "
$ws.Range("D3").Value = "2.284.86"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "501.79"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("D6").Value = "130.01"
$ws.Range("E6").Value = "  +2.18%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  +1.08%  "
$ws.Range("E10").Value = "  +0.88%  "
$ws.Range("D11").Value = "0.339"
$ws.Range("E11").Value = "  +4.95%  "
$ws.Range("E12").Value = "  +2.34%  "
$ws.Range("D13").Value = "2.693.38"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").Value = "23.06"
$ws.Range("E14").Value = "  +6.56%  "
$ws.Range("D15").Value = "54.392.14"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "2.292.71"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "10.28"
$ws.Range("E18").Value = "  +2.64%  "
$ws.Range("D19").Value = "4.19"
$ws.Range("E19").Value = "  +3.07%  "
$ws.Range("D20").Value = "304.00"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "6.35"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "61.88"
$ws.Range("E23").Value = "  -2.30%  "
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").Value = "7.35"
$ws.Range("E26").Value = "  +3.65%  "
$ws.Range("D27").Value = "170.69"
$ws.Range("E27").Value = "  +1.21%  "
$ws.Range("E28").Value = "  +1.72%  "
$ws.Range("D29").Value = "0.0₃0692"
$ws.Range("E29").Value = "  +1.06%  "
$ws.Range("D30").Value = "5.94"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("D33").Value = "17.80"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").Value = "0.955"
$ws.Range("E34").Value = "  +10.12%  "
$ws.Range("E35").Value = "  -0.16%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "3.74"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  +1.30%  "
$ws.Range("E40").Value = "  +1.28%  "
$ws.Range("D41").Value = "4.88"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").Value = "125.98"
$ws.Range("E43").Value = "  +3.17%  "
$ws.Range("E44").Value = "  +0.67%  "
$ws.Range("D45").Value = "0.549"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("D46").Value = "242.66"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  +1.71%  "
$ws.Range("D49").Value = "10.78"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "16.40"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("E51").Value = "  -0.50%  "
